$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N112").Value = -7237.0871
$ws.Range("J112").Value = 1673.6957
$ws.Range("L112").Value = 5021.0871
$ws.Range("H112").Value = 1673.6957
$ws.Range("N138").Value = -32367.6362
$ws.Range("J138").Value = 7362.5454
$ws.Range("L138").Value = 22087.6362
$ws.Range("H138").Value = 6075.1113
$ws.Range("I141").Value = 1452.1904
$ws.Range("N141").Value = -40402
$ws.Range("M141").Value = 823.4287999999997
$ws.Range("J141").Value = 10014
$ws.Range("L141").Value = 30042
$ws.Range("H141").Value = 3354.8147
$ws.Range("K141").Value = 4356.5712

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I45").Value = 21682.2
$ws.Range("N45").Value = -3095.25
$ws.Range("M45").Value = -21305.2
$ws.Range("J45").Value = 2341.25
$ws.Range("L45").Value = 2341.25
$ws.Range("H45").Value = 8029.7646
$ws.Range("K45").Value = 21682.2
$ws.Range("I74").Value = 49248.094
$ws.Range("M74").Value = -48374.094
$ws.Range("H74").Value = 36785.137
$ws.Range("K74").Value = 49248.094
$ws.Range("I77").Value = 49248.094
$ws.Range("M77").Value = -241872.47
$ws.Range("H77").Value = 36785.137
$ws.Range("K77").Value = 246240.47
$ws.Range("I102").Value = 13375.875
$ws.Range("N102").Value = -25696.2
$ws.Range("M102").Value = -11753.875
$ws.Range("J102").Value = 22452.2
$ws.Range("L102").Value = 22452.2
$ws.Range("H102").Value = 16866.77
$ws.Range("K102").Value = 13375.875
$ws.Range("I122").Value = 3553.1667
$ws.Range("N122").Value = -19000
$ws.Range("M122").Value = -8209.500100000001
$ws.Range("J122").Value = 4700
$ws.Range("L122").Value = 14100
$ws.Range("H122").Value = 3839.875
$ws.Range("K122").Value = 10659.5001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I20").Value = 144733.47
$ws.Range("M20").Value = -144486.47
$ws.Range("H20").Value = 144733.47
$ws.Range("K20").Value = 144733.47
$ws.Range("I36").Value = 680.75
$ws.Range("N36").Value = -1398
$ws.Range("M36").Value = -146.75
$ws.Range("J36").Value = 330
$ws.Range("L36").Value = 330
$ws.Range("H36").Value = 610.6
$ws.Range("K36").Value = 680.75
$ws.Range("I94").Value = 952.6070999999999
$ws.Range("N94").Value = -3494.5715
$ws.Range("M94").Value = -501.6070999999999
$ws.Range("J94").Value = 2592.5715
$ws.Range("L94").Value = 2592.5715
$ws.Range("H94").Value = 1499.262
$ws.Range("K94").Value = 952.6070999999999
$ws.Range("I99").Value = 2966.6667
$ws.Range("M99").Value = -1468.6667
$ws.Range("H99").Value = 2407456.2
$ws.Range("K99").Value = 2966.6667
$ws.Range("I105").Value = 1381.0555
$ws.Range("N105").Value = -6106.25
$ws.Range("M105").Value = 365.9445000000001
$ws.Range("J105").Value = 2612.25
$ws.Range("L105").Value = 2612.25
$ws.Range("H105").Value = 1604.909
$ws.Range("K105").Value = 1381.0555
$ws.Range("I107").Value = 4946.724
$ws.Range("M107").Value = -3026.724
$ws.Range("H107").Value = 5129.706
$ws.Range("K107").Value = 4946.724

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I4").Value = 635
$ws.Range("N4").Value = -7874
$ws.Range("M4").Value = -523
$ws.Range("J4").Value = 7650
$ws.Range("L4").Value = 7650
$ws.Range("H4").Value = 2639.2856
$ws.Range("K4").Value = 635
$ws.Range("I7").Value = 9138.416999999999
$ws.Range("N7").Value = -10551.6
$ws.Range("M7").Value = -9025.416999999999
$ws.Range("J7").Value = 10325.6
$ws.Range("L7").Value = 10325.6
$ws.Range("H7").Value = 9678.046
$ws.Range("K7").Value = 9138.416999999999
$ws.Range("I13").Value = 16
$ws.Range("N13").Value = -50279
$ws.Range("M13").Value = 123
$ws.Range("J13").Value = 50001
$ws.Range("L13").Value = 50001
$ws.Range("H13").Value = 25008.5
$ws.Range("K13").Value = 16
$ws.Range("I31").Value = 3182.1562
$ws.Range("N31").Value = -7784.8716
$ws.Range("M31").Value = -2887.1562
$ws.Range("J31").Value = 7194.8716
$ws.Range("L31").Value = 7194.8716
$ws.Range("H31").Value = 5386.3237
$ws.Range("K31").Value = 3182.1562
$ws.Range("I34").Value = 3182.1562
$ws.Range("N34").Value = -7598.8716
$ws.Range("M34").Value = -2980.1562
$ws.Range("J34").Value = 7194.8716
$ws.Range("L34").Value = 7194.8716
$ws.Range("H34").Value = 5386.3237
$ws.Range("K34").Value = 3182.1562
$ws.Range("I99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H99").Value = 15627612
$ws.Range("K99").Value = 0
$ws.Range("I105").Value = 1907.9231
$ws.Range("M105").Value = -160.9231
$ws.Range("H105").Value = 2550.8125
$ws.Range("K105").Value = 1907.9231
$ws.Range("N116").Value = -87543.875
$ws.Range("J116").Value = 78365.875
$ws.Range("L116").Value = 78365.875
$ws.Range("H116").Value = 78365.875
$ws.Range("I126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H126").Value = 15627612
$ws.Range("K126").Value = 0
$ws.Range("I132").Value = 2438.75
$ws.Range("M132").Value = -4786.25
$ws.Range("H132").Value = 3884
$ws.Range("K132").Value = 7316.25
$ws.Range("N141").Value = -359345.47
$ws.Range("J141").Value = 348985.47
$ws.Range("L141").Value = 348985.47
$ws.Range("H141").Value = 348985.47

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I31").Value = 0
$ws.Range("N31").Value = -33473.001
$ws.Range("M31").ClearContents()
$ws.Range("J31").Value = 10965.667
$ws.Range("L31").Value = 32897.001
$ws.Range("H31").Value = 10965.667
$ws.Range("K31").Value = 0
$ws.Range("I81").Value = 895.34784
$ws.Range("M81").Value = -1563.04352
$ws.Range("H81").Value = 4635.5713
$ws.Range("K81").Value = 2686.04352
$ws.Range("I84").Value = 895.34784
$ws.Range("M84").Value = -2442.130560000001
$ws.Range("H84").Value = 4635.5713
$ws.Range("K84").Value = 8058.130560000001
$ws.Range("I107").Value = 276.44446
$ws.Range("N107").Value = -5729.4999
$ws.Range("M107").Value = 1090.66662
$ws.Range("J107").Value = 629.8333
$ws.Range("L107").Value = 1889.4999
$ws.Range("H107").Value = 417.8
$ws.Range("K107").Value = 829.33338
$ws.Range("I113").Value = 2901.5
$ws.Range("M113").Value = -6534.5
$ws.Range("H113").Value = 3575.75
$ws.Range("K113").Value = 8704.5
$ws.Range("I129").Value = 908
$ws.Range("M129").Value = 2276
$ws.Range("H129").Value = 1216.8
$ws.Range("K129").Value = 2724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I70").Value = 107344.4
$ws.Range("M70").Value = -107074.4
$ws.Range("H70").Value = 75346.13
$ws.Range("K70").Value = 107344.4
$ws.Range("I73").Value = 107344.4
$ws.Range("M73").Value = -106408.4
$ws.Range("H73").Value = 75346.13
$ws.Range("K73").Value = 107344.4
$ws.Range("N96").Value = -27912.334
$ws.Range("J96").Value = 22420.334
$ws.Range("L96").Value = 22420.334
$ws.Range("H96").Value = 22420.334
$ws.Range("N97").Value = -1988.5
$ws.Range("J97").Value = 996.5
$ws.Range("L97").Value = 996.5
$ws.Range("H97").Value = 782.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I18").Value = 125
$ws.Range("N18").ClearContents()
$ws.Range("M18").Value = 47
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("H18").Value = 125
$ws.Range("K18").Value = 125
$ws.Range("I20").Value = 150
$ws.Range("N20").Value = -1434766.2
$ws.Range("M20").Value = 76
$ws.Range("J20").Value = 1434314.2
$ws.Range("L20").Value = 1434314.2
$ws.Range("H20").Value = 912800
$ws.Range("K20").Value = 150
$ws.Range("I55").Value = 692
$ws.Range("N55").Value = -2321.1
$ws.Range("M55").Value = -519
$ws.Range("J55").Value = 1975.1
$ws.Range("L55").Value = 1975.1
$ws.Range("H55").Value = 1185.5
$ws.Range("K55").Value = 692
$ws.Range("N100").Value = -4582
$ws.Range("J100").Value = 3500
$ws.Range("L100").Value = 3500
$ws.Range("H100").Value = 35500
$ws.Range("I132").Value = 3478.5789
$ws.Range("M132").Value = -7905.736699999999
$ws.Range("H132").Value = 3922.6785
$ws.Range("K132").Value = 10435.7367
$ws.Range("I136").Value = 7706.533
$ws.Range("M136").Value = -20569.599
$ws.Range("H136").Value = 6554.1724
$ws.Range("K136").Value = 23119.599
$ws.Range("N137").Value = -88137.5
$ws.Range("J137").Value = 77937.5
$ws.Range("L137").Value = 77937.5
$ws.Range("H137").Value = 72856.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I81").Value = 1195.1538
$ws.Range("N81").Value = -69002
$ws.Range("M81").Value = -1329.3076
$ws.Range("J81").Value = 33440
$ws.Range("L81").Value = 66880
$ws.Range("H81").Value = 7241.0625
$ws.Range("K81").Value = 2390.3076
$ws.Range("I84").Value = 1195.1538
$ws.Range("N84").Value = -345008
$ws.Range("M84").Value = -6647.538
$ws.Range("J84").Value = 33440
$ws.Range("L84").Value = 334400
$ws.Range("H84").Value = 7241.0625
$ws.Range("K84").Value = 11951.538
$ws.Range("I96").Value = 1504.9
$ws.Range("M96").Value = -131.9000000000001
$ws.Range("H96").Value = 1648.3684
$ws.Range("K96").Value = 1504.9
$ws.Range("N109").Value = -69422
$ws.Range("J109").Value = 66648
$ws.Range("L109").Value = 66648
$ws.Range("H109").Value = 64432
$ws.Range("I126").Value = 1773.24
$ws.Range("N126").Value = -158315.75
$ws.Range("M126").Value = -2849.72
$ws.Range("J126").Value = 51125.25
$ws.Range("L126").Value = 153375.75
$ws.Range("H126").Value = 8580.414000000001
$ws.Range("K126").Value = 5319.72
